# Two new bank-statement rows were recorded at the top of the sheet
# (most-recent-first ledger), pushing the existing 30 transaction rows
# down by two, and two more blank trailer rows were appended at the
# bottom. The H1 formula (which mirrors row 1 into a PHP array-literal
# string) is refreshed to match the new row 1 and to stamp the
# "mo_fecha_crea" field with NOW() instead of a fixed literal, plus a
# new mo_borrado_logico field.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Detach the H1 formula first so inserting rows above it does not drag
# it down into H3 (the refreshed formula belongs back on row 1 only).
$ws.Range("H1").ClearContents()

# Shift the existing 30 data rows (1-30) down to rows 3-32.
$ws.Rows("1:2").Insert()

# Append two blank trailer rows at the bottom (was 33-34, now 33-36).
$ws.Range("A35").NumberFormat = "m/d/yy"
$ws.Range("A36").NumberFormat = "m/d/yy"

# Fill in the two new transactions at the top.
$ws.Range("A1").Value = 41681
$ws.Range("A1").NumberFormat = "m/d/yy"
$ws.Range("B1").Value = 'INTERES A SU FAVOR'
$ws.Range("C1").Value = 'C'
$ws.Range("D1").Value = '0000950673'
$ws.Range("E1").Value = 'AGENCIA PARA PROCESOS BATCH'
$ws.Range("F1").Value = '0.25  '
$ws.Range("G1").Value = '4014.30'

$ws.Range("A2").Value = 41680
$ws.Range("A2").NumberFormat = "m/d/yy"
$ws.Range("B2").Value = 'INTERES A SU FAVOR'
$ws.Range("C2").Value = 'C'
$ws.Range("D2").Value = '0000950683'
$ws.Range("E2").Value = 'AGENCIA PARA PROCESOS BATCH'
$ws.Range("F2").Value = '0.25  '
$ws.Range("G2").Value = '4014.05'

# Refresh the H1 formula to reference the new row-1 data and to use
# NOW() for mo_fecha_crea, with the new mo_borrado_logico field.
$ws.Range("H1").Formula = '=CONCATENATE("array(''mo_fecha'' => new \DateTime(''",TEXT(A1,"yyyy-mm-dd"),"''), ''mo_concepto'' => ''",B1,"'', ''mo_tipo'' => ''",C1,"'', ''mo_documento'' => ''",D1,"'', ''mo_oficina'' => ''",E1,"'', ''mo_monto'' => ",F1,", ''mo_saldo'' => ",G1,", ''mo_fecha_crea'' => new \DateTime(''",TEXT(NOW(),"yyyy-mm-dd H:m:s"),"''), ''mo_quien_crea'' => 1, ''mo_fecha_modifica'' => NULL, ''mo_quien_modifica'' => NULL ''mo_borrado_logico'' => false),")'

# Restore the active-cell selection shown in the saved workbook.
$ws.Range("H13").Select()
